# "Change 14 bus test"
#
# Loosens the reactive-power limits (Qmin/Qmax) on several generator buses
# and tweaks their scheduled real power, relaxes three transformer turns
# ratios to 1 on the IEEE network-line sheet, appends eight new
# (near-)infinite-impedance self/loop branch rows to the same sheet, and
# disables the "plot admittance" option on the Advance sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Bus sheet: PGi (E), Qmin (I), Qmax (J) for buses 2, 3, 6, 8
# ---------------------------------------------------------------------------
$wsBus = $wb.Worksheets.Item("Bus")

$wsBus.Range("E4").Value = 0.2
$wsBus.Range("I4").Value = -999
$wsBus.Range("J4").Value = 999

$wsBus.Range("E5").Value = 0.1
$wsBus.Range("I5").Value = -999
$wsBus.Range("J5").Value = 999

$wsBus.Range("E8").Value = 0.1
$wsBus.Range("I8").Value = -999
$wsBus.Range("J8").Value = 999

$wsBus.Range("E10").Value = 0.1
$wsBus.Range("I10").Value = -999
$wsBus.Range("J10").Value = 999

# ---------------------------------------------------------------------------
# NetworkLine_IEEE sheet: a few transformer turns ratios -> 1
# ---------------------------------------------------------------------------
$wsNLI = $wb.Worksheets.Item("NetworkLine_IEEE")

$wsNLI.Range("G13").Value = 1
$wsNLI.Range("G14").Value = 1
$wsNLI.Range("G15").Value = 1

# Append eight new self/loop branch rows (27-34): From bus / To bus / R / X
# ("inf") / B (1E-3, scientific fmt) / status enable (0) / turns ratio (1)
$newBranchBuses = @(6, 7, 8, 10, 11, 12, 13, 14)
$row = 27
foreach ($bus in $newBranchBuses) {
    $wsNLI.Cells.Item($row, 1).Value = $bus
    $wsNLI.Cells.Item($row, 2).Value = $bus
    $wsNLI.Cells.Item($row, 3).Value = "inf"
    $wsNLI.Cells.Item($row, 4).Value = "inf"
    $wsNLI.Cells.Item($row, 5).Value = 0.001
    $wsNLI.Cells.Item($row, 5).NumberFormat = "0.00E+00"
    $wsNLI.Cells.Item($row, 6).Value = 0
    $wsNLI.Cells.Item($row, 7).Value = 1
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# Advance sheet: disable "Enable (plot admittance)"
# ---------------------------------------------------------------------------
$wsAdv = $wb.Worksheets.Item("Advance")
$wsAdv.Range("B10").Value = 0

# ---------------------------------------------------------------------------
# Selections / active sheet. Select in an order such that "Bus" ends up the
# last-activated (and therefore tabSelected) sheet.
# ---------------------------------------------------------------------------
$wsApp = $wb.Worksheets.Item("Apparatus")
[void]$wsApp.Range("E3").Select()

[void]$wsNLI.Range("E30").Select()

[void]$wsAdv.Range("B11").Select()

[void]$wsBus.Range("E5").Select()
